{"js": "const body = context.document.body;\nconst results = body.search(\"Version 2.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Version 1.\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Version 2.\"\n$newText = \"Version 1.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif ($found) {\n    $matchRange = $find.Parent\n    $matchText = $matchRange.Text\n\n    # Replace only the characters that actually differ so in-between\n    # markup (bookmarks, proofing marks, etc.) inside the matched range\n    # is left untouched.\n    $prefixLen = 0\n    $minLen = [Math]::Min($matchText.Length, $newText.Length)\n    while ($prefixLen -lt $minLen -and $matchText[$prefixLen] -eq $newText[$prefixLen]) {\n        $prefixLen++\n    }\n\n    $suffixLen = 0\n    while ($suffixLen -lt ($minLen - $prefixLen) -and\n           $matchText[$matchText.Length - 1 - $suffixLen] -eq $newText[$newText.Length - 1 - $suffixLen]) {\n        $suffixLen++\n    }\n\n    $subStart = $matchRange.Start + $prefixLen\n    $subEnd = $matchRange.End - $suffixLen\n    $replacement = $newText.Substring($prefixLen, $newText.Length - $prefixLen - $suffixLen)\n\n    $sub = $d.Range($subStart, $subEnd)\n    $sub.Text = $replacement\n}\n"}
